$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New matchup rows for "spring 23 week 10 inputs".
# Existing data occupies rows 2..1758 (row 1 = headers); append 24 new rows
# starting at row 1759 through row 1782.
$rows = @(
    @(4, 0, 4, 2),
    @(5, 0, 7, 3),
    @(4, 1, 3, 2),
    @(6, 1, 5, 2),
    @(3, 2, 2, 1),
    @(4, 0, 3, 3),
    @(3, 1, 3, 2),
    @(6, 1, 4, 2),
    @(4, 2, 6, 0),
    @(6, 1, 3, 2),
    @(4, 2, 4, 0),
    @(5, 3, 5, 0),
    @(4, 0, 3, 2),
    @(6, 1, 5, 2),
    @(3, 3, 3, 0),
    @(5, 2, 5, 0),
    @(3, 2, 3, 1),
    @(4, 0, 4, 3),
    @(7, 0, 7, 2),
    @(4, 1, 6, 2),
    @(3, 0, 7, 3),
    @(5, 2, 3, 1),
    @(4, 0, 2, 2),
    @(5, 2, 5, 1)
)

$startRow = 1759
$numRows = $rows.Count
$numCols = 4

$arr = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}

$endRow = $startRow + $numRows - 1
$ws.Range("A$startRow`:D$endRow").Value = $arr

# Match the viewport/selection state recorded in the committed workbook.
$excel.ActiveWindow.ScrollRow = 1770
$ws.Range("A$($endRow + 1)").Select()
